# Applies the "Avoid the word you(r) in Readme" / JJ.Framework reference fix
# edits described in the commit diff.
#
# Strategy: most of the diff is just merging several adjacent <w:r> runs
# (that share identical formatting, often separated by <w:proofErr/> spell
# check markers) back into a single run, with a couple of small text
# tweaks, plus one structural change (splitting a checklist item into two,
# and relocating the "_GoBack" bookmark into the new second item).
#
# Because Range.Text assignment cannot precisely control run/bookmark
# boundaries, we use Range.InsertXML with hand-built OOXML fragments for
# full control, matching exactly what the target XML should look like.
#
# We edit paragraphs from the bottom of the document upwards so that
# paragraph indices we have already looked up stay valid even after an
# edit inserts an extra paragraph earlier in the flow.

function Set-ParaXml {
    param($doc, $paraIndex, $innerXml)
    $p = $doc.Paragraphs($paraIndex)
    $r = $p.Range
    # Exclude the trailing paragraph mark so we only replace the paragraph's
    # content (and keep it as a single paragraph unless innerXml itself
    # contains multiple <w:p> elements, in which case Word splits it).
    $body = $doc.Range($r.Start, $r.End - 1)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
          '<pkg:xmlData>' + `
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
              '<w:body>' + $innerXml + '</w:body>' + `
            '</w:document>' + `
          '</pkg:xmlData>' + `
        '</pkg:part>' + `
      '</pkg:package>'
    $body.InsertXML($pkg)
}

$d = $word.ActiveDocument

# --- Paragraph 72: "I turned out to add all the JJ.Framework csproj's ..." ---
# merge runs, drop proofErr spell-check wrappers
$xml72 = '<w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="25"/></w:numPr><w:rPr><w:color w:val="BFBFBF"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:color w:val="BFBFBF"/><w:lang w:val="en-US"/></w:rPr><w:t>I turned out to add all the JJ.Framework csproj' + [char]0x2019 + 's that the dependencies asked for, because I forgot a few. That fixes it for my project.</w:t></w:r>' + `
  '</w:p>'
Set-ParaXml $d 72 $xml72

# --- Paragraph 71: "A JJ.Framework.WinForms.TestForms Form will also open in de designer." ---
$xml71 = '<w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="25"/></w:numPr><w:rPr><w:color w:val="BFBFBF"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:color w:val="BFBFBF"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">A JJ.Framework.WinForms.TestForms Form will also open in de designer.</w:t></w:r>' + `
  '</w:p>'
Set-ParaXml $d 71 $xml71

# --- Paragraph 70: "The JJ.Framework.WinForms.TestForms seems to run fine." ---
$xml70 = '<w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="25"/></w:numPr><w:rPr><w:color w:val="BFBFBF"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:color w:val="BFBFBF"/><w:lang w:val="en-US"/></w:rPr><w:t>The JJ.Framework.WinForms.TestForms seems to run fine.</w:t></w:r>' + `
  '</w:p>'
Set-ParaXml $d 70 $xml70

# --- Paragraph 68: "[x] Error placing DiagramControl on Form: cannot load JJ.Framework.VectorGraphics." ---
$xml68 = '<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="25"/></w:numPr><w:rPr><w:color w:val="BFBFBF"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:color w:val="BFBFBF"/><w:lang w:val="en-US"/></w:rPr><w:t>[x] Error placing DiagramControl on Form: cannot load JJ.Framework.VectorGraphics.</w:t></w:r>' + `
  '</w:p>'
Set-ParaXml $d 68 $xml68

# --- Paragraph 67: "ElementPosition.SetMarginInPixels() based on code from CurveDetailsViewModelToDiagramConverter around line 148?" ---
# keep the two original runs (different rsid markers in source), just drop proofErr wrappers
$xml67 = '<w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="25"/></w:numPr><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr><w:t>ElementPosition.</w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr><w:t>SetMarginInPixels</w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">() based on code from </w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr><w:t>CurveDetailsViewModelToDiagramConverter</w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> around line 148?</w:t></w:r>' + `
  '</w:p>'
Set-ParaXml $d 67 $xml67

# --- Paragraph 63: "[x] Gave a VectorGraphics Element.Children a Clear method." ---
$xml63 = '<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="25"/></w:numPr><w:rPr><w:color w:val="BFBFBF"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:color w:val="BFBFBF"/><w:lang w:val="en-US"/></w:rPr><w:t>[</w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="BFBFBF"/><w:lang w:val="en-US"/></w:rPr><w:t>x</w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="BFBFBF"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">] </w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="BFBFBF"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Gave </w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="BFBFBF"/><w:lang w:val="en-US"/></w:rPr><w:t>a VectorGraphics Element.Children a Clear method.</w:t></w:r>' + `
  '</w:p>'
Set-ParaXml $d 63 $xml63

# --- Paragraph 59: "2020-08-05 MedsUseInfoGraphic Details" ---
$xml59 = '<w:p><w:pPr><w:pStyle w:val="Heading3"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>2020-08-05 MedsUseInfoGraphic Details</w:t></w:r>' + `
  '</w:p>'
Set-ParaXml $d 59 $xml59

# --- Paragraph 56: "[ ] JJ.MedsUseInfographic.Data.SpecialFormat or .FromNotes: Parsing text from how I type in my meds use in my notes." ---
$xml56 = '<w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="25"/></w:numPr><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr><w:t>[ ]</w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> JJ.MedsUseInfographic.Data.SpecialFormat</w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> or .FromNotes</w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr><w:t>:</w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Parsing </w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">text </w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">from </w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr><w:t>h</w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr><w:t>ow I type in my meds use in my notes.</w:t></w:r>' + `
  '</w:p>'
Set-ParaXml $d 56 $xml56

# --- Paragraph 36: "[ ] NuGet packaging JJ.Framework.WinForms?" ---
$xml36 = '<w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="25"/></w:numPr><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:color w:val="FFC000"/><w:lang w:val="en-US"/></w:rPr><w:t>[ ] NuGet packaging JJ.Framework.WinForms?</w:t></w:r>' + `
  '</w:p>'
Set-ParaXml $d 36 $xml36

# --- Paragraph 34: split into two checklist items; relocate the "_GoBack" bookmark here ---
$xml34 = '<w:p><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="25"/></w:numPr><w:rPr><w:color w:val="92D050"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:color w:val="92D050"/><w:lang w:val="en-US"/></w:rPr><w:t>[</w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="92D050"/><w:lang w:val="en-US"/></w:rPr><w:t>x</w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="92D050"/><w:lang w:val="en-US"/></w:rPr><w:t>] Alternative: Move contents of JJs Software Small folder to the JJs Software folder?</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="25"/></w:numPr><w:rPr><w:color w:val="92D050"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:color w:val="92D050"/><w:lang w:val="en-US"/></w:rPr><w:t>[ ] JJ.Framework references are still not found.</w:t></w:r>' + `
  '<w:bookmarkStart w:id="4" w:name="_GoBack"/><w:bookmarkEnd w:id="4"/>' + `
  '</w:p>'
Set-ParaXml $d 34 $xml34

# --- Paragraph 30: merge runs, drop the (now relocated) "_GoBack" bookmark ---
$xml30 = '<w:p><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="25"/></w:numPr><w:rPr><w:color w:val="92D050"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:color w:val="92D050"/><w:lang w:val="en-US"/></w:rPr><w:t>[ ] ..</w:t></w:r>' + `
  '<w:r><w:rPr><w:color w:val="92D050"/><w:lang w:val="en-US"/></w:rPr><w:t>JJ.Framework is not in the right folder</w:t></w:r>' + `
  '</w:p>'
Set-ParaXml $d 30 $xml30

Write-Host "Edits applied."
